$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.022.32'
$ws.Range('E2').Value = '  +2.00%  '
$ws.Range('D3').Value = '2.299.78'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '309.87'
$ws.Range('E5').Value = '  +1.64%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '100.78'
$ws.Range('E6').Value = '  +4.33%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.530'
$ws.Range('E7').Value = '  +0.13%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.506'
$ws.Range('E9').Value = '  +2.95%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.09'
$ws.Range('E10').Value = '  +1.52%  '
$ws.Range('E11').Value = '  +2.99%  '
$ws.Range('E12').Value = '  +0.72%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.95'
$ws.Range('E13').Value = '  +4.60%  '
$ws.Range('D14').Value = '2.656.08'
$ws.Range('E14').Value = '  +1.40%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.83'
$ws.Range('E15').Value = '  +2.67%  '
$ws.Range('D16').Value = '2.309.70'
$ws.Range('E16').Value = '  +1.57%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.802'
$ws.Range('E17').Value = '  +0.95%  '
$ws.Range('D18').Value = '43.006.34'
$ws.Range('E18').Value = '  +2.18%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.51'
$ws.Range('E19').Value = '  +0.46%  '
$ws.Range('E20').Value = '  +0.90%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.07'
$ws.Range('E21').Value = '  +1.17%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.08'
$ws.Range('E22').Value = '  +0.34%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '239.83'
$ws.Range('E23').Value = '  +0.59%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.03'
$ws.Range('E24').Value = '  +4.26%  '
$ws.Range('E25').Value = '  +1.50%  '
$ws.Range('E26').Value = '  +0.14%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '23.94'
$ws.Range('E27').Value = '  +0.90%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '38.31'
$ws.Range('E28').Value = '  +2.82%  '
$ws.Range('B29').Value = 'Cosmos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.63'
$ws.Range('E29').Value = '  +1.04%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.12'
$ws.Range('E30').Value = '  +0.87%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '165.47'
$ws.Range('E31').Value = '  +3.77%  '
$ws.Range('E32').Value = '  +0.89%  '
$ws.Range('E33').Value = '  +0.09%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.15'
$ws.Range('E34').Value = '  -1.27%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '17.73'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0739'
$ws.Range('E36').Value = '  -0.06%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.38'
$ws.Range('E37').Value = '  +0.62%  '
$ws.Range('E38').Value = '  -0.01%  '
$ws.Range('E39').Value = '  +0.28%  '
$ws.Range('E40').Value = '  +1.27%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.18'
$ws.Range('E41').Value = '  +2.33%  '
$ws.Range('E42').Value = '  -5.59%  '
$ws.Range('E43').Value = '  +1.61%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '19.23'
$ws.Range('E44').Value = '  +0.78%  '
$ws.Range('D45').Value = '1.963.81'
$ws.Range('E45').Value = '  -1.32%  '
$ws.Range('E46').Value = '  +3.20%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.84'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.02'
$ws.Range('E48').Value = '  +21.06%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '54.85'
$ws.Range('E49').Value = '  +3.16%  '
$ws.Range('D50').Value = '2.525.12'
$ws.Range('E50').Value = '  +1.41%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.54'
$ws.Range('E51').Value = '  +2.01%  '
